$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Stable style-source cell (untouched by this edit, same style index as the
# "General"-formatted text cells we are about to update) used to restore the
# original cell style after writing values that Excel would otherwise
# auto-convert (e.g. trailing "%" reinterpreted as a percentage number).
$styleSource = $ws.Range("C2")

$ws.Range("E2").Value = "2026-02-16 17:18:37"
$ws.Range("I2").Value = "17.4 mm"
$ws.Range("K2").Value = "4.2 MJ/m2"
$ws.Range("E3").Value = "2026-02-16 17:18:40"
$ws.Range("I3").Value = "9.4 mm"
$ws.Range("K3").Value = "6.7 MJ/m2"
$ws.Range("N3").Value = "-1.9 °C 16:55 TU"
$ws.Range("E4").Value = "2026-02-16 17:18:42"
$styleSource.Copy() | Out-Null
$ws.Range("H4").NumberFormat = "@"
$ws.Range("H4").Value = "58%"
$ws.Range("H4").PasteSpecial(-4122) | Out-Null
$ws.Range("J4").Value = "1012.5 hPa"
$ws.Range("K4").Value = "13.5 MJ/m2"
$ws.Range("O4").Value = "14.0 °C"
$ws.Range("E5").Value = "2026-02-16 17:18:45"
$ws.Range("G5").Value = "144 cm"
$ws.Range("I5").Value = "20.1 mm"
$ws.Range("L5").Value = "35.6 km/h - 340º 16:46 TU"
$ws.Range("N5").Value = "-1.2 °C 16:59 TU"
$ws.Range("E6").Value = "2026-02-16 17:18:47"
$ws.Range("J6").Value = "1012.6 hPa"
$ws.Range("K6").Value = "13.2 MJ/m2"
$ws.Range("O6").Value = "11.5 °C"
$ws.Range("E7").Value = "2026-02-16 17:18:50"
$ws.Range("J7").Value = "1013.6 hPa"
$ws.Range("K7").Value = "11.5 MJ/m2"
$ws.Range("O7").Value = "16.1 °C"
$ws.Range("E8").Value = "2026-02-16 17:18:52"
$ws.Range("J8").Value = "1013.2 hPa"
$ws.Range("O8").Value = "12.3 °C"
$ws.Range("E9").Value = "2026-02-16 17:18:55"
$styleSource.Copy() | Out-Null
$ws.Range("H9").NumberFormat = "@"
$ws.Range("H9").Value = "74%"
$ws.Range("H9").PasteSpecial(-4122) | Out-Null
$ws.Range("O9").Value = "10.9 °C"
$ws.Range("E10").Value = "2026-02-16 17:18:57"
$ws.Range("O10").Value = "10.7 °C"
$ws.Range("E11").Value = "2026-02-16 17:18:59"
$ws.Range("O11").Value = "6.3 °C"
$ws.Range("E12").Value = "2026-02-16 17:19:02"
$ws.Range("O12").Value = "10.3 °C"
$ws.Range("E13").Value = "2026-02-16 17:19:04"
$ws.Range("J13").Value = "1015.2 hPa"
$ws.Range("K13").Value = "12.7 MJ/m2"
$ws.Range("O13").Value = "5.3 °C"
$ws.Range("E14").Value = "2026-02-16 17:19:07"
$ws.Range("K14").Value = "9.9 MJ/m2"
$ws.Range("O14").Value = "15.8 °C"
$ws.Range("E15").Value = "2026-02-16 17:19:09"
$ws.Range("O15").Value = "11.0 °C"
$ws.Range("E16").Value = "2026-02-16 17:19:11"
$ws.Range("K16").Value = "9.5 MJ/m2"
$ws.Range("E17").Value = "2026-02-16 17:19:13"
$styleSource.Copy() | Out-Null
$ws.Range("H17").NumberFormat = "@"
$ws.Range("H17").Value = "67%"
$ws.Range("H17").PasteSpecial(-4122) | Out-Null
$ws.Range("K17").Value = "14.1 MJ/m2"
$ws.Range("N17").Value = "4.7 °C 16:59 TU"
$ws.Range("E18").Value = "2026-02-16 17:19:16"
$ws.Range("K18").Value = "13.6 MJ/m2"
$ws.Range("L18").Value = "26.3 km/h - 36º 16:49 TU"
$ws.Range("O18").Value = "10.7 °C"
$ws.Range("E19").Value = "2026-02-16 17:19:18"
$styleSource.Copy() | Out-Null
$ws.Range("H19").NumberFormat = "@"
$ws.Range("H19").Value = "83%"
$ws.Range("H19").PasteSpecial(-4122) | Out-Null
$ws.Range("K19").Value = "13.0 MJ/m2"
$ws.Range("O19").Value = "6.6 °C"
$ws.Range("E20").Value = "2026-02-16 17:19:21"
$ws.Range("K20").Value = "10.9 MJ/m2"
$ws.Range("E21").Value = "2026-02-16 17:19:23"
$ws.Range("J21").Value = "1014.6 hPa"
$ws.Range("O21").Value = "8.1 °C"
$ws.Range("E22").Value = "2026-02-16 17:19:26"
$ws.Range("E23").Value = "2026-02-16 17:19:28"
$ws.Range("I23").Value = "11.9 mm"
$ws.Range("K23").Value = "8.7 MJ/m2"
$ws.Range("N23").Value = "-1.8 °C 16:55 TU"
$ws.Range("E24").Value = "2026-02-16 17:19:31"
$ws.Range("K24").Value = "11.6 MJ/m2"
$ws.Range("O24").Value = "12.8 °C"
$ws.Range("E25").Value = "2026-02-16 17:19:33"
$ws.Range("I25").Value = "5.1 mm"
$ws.Range("K25").Value = "6.7 MJ/m2"
$ws.Range("E26").Value = "2026-02-16 17:19:36"
$ws.Range("E27").Value = "2026-02-16 17:19:38"
$ws.Range("K27").Value = "7.8 MJ/m2"
$ws.Range("E28").Value = "2026-02-16 17:19:40"
$ws.Range("J28").Value = "1013.1 hPa"
$ws.Range("K28").Value = "12.9 MJ/m2"
$ws.Range("O28").Value = "9.2 °C"
$ws.Range("E29").Value = "2026-02-16 17:19:43"
$ws.Range("O29").Value = "10.6 °C"
$ws.Range("E30").Value = "2026-02-16 17:19:45"
$ws.Range("K30").Value = "12.3 MJ/m2"
$ws.Range("O30").Value = "11.8 °C"
$ws.Range("E31").Value = "2026-02-16 17:19:48"
$ws.Range("E32").Value = "2026-02-16 17:19:50"
$styleSource.Copy() | Out-Null
$ws.Range("H32").NumberFormat = "@"
$ws.Range("H32").Value = "79%"
$ws.Range("H32").PasteSpecial(-4122) | Out-Null
$ws.Range("K32").Value = "11.1 MJ/m2"
$ws.Range("O32").Value = "8.4 °C"
$ws.Range("E33").Value = "2026-02-16 17:19:53"
$ws.Range("J33").Value = "1014.0 hPa"
$ws.Range("K33").Value = "9.9 MJ/m2"
$ws.Range("O33").Value = "6.2 °C"
$ws.Range("E34").Value = "2026-02-16 17:19:55"
$ws.Range("K34").Value = "8.8 MJ/m2"
$ws.Range("E35").Value = "2026-02-16 17:19:58"
$ws.Range("J35").Value = "1017.0 hPa"
$ws.Range("L35").Value = "79.9 km/h - 268º 16:57 TU"
$ws.Range("O35").Value = "9.4 °C"
$ws.Range("E36").Value = "2026-02-16 17:20:00"
$ws.Range("O36").Value = "11.4 °C"
$ws.Range("E37").Value = "2026-02-16 17:20:03"
$styleSource.Copy() | Out-Null
$ws.Range("H37").NumberFormat = "@"
$ws.Range("H37").Value = "82%"
$ws.Range("H37").PasteSpecial(-4122) | Out-Null
$ws.Range("J37").Value = "1015.1 hPa"
$ws.Range("O37").Value = "6.2 °C"
$ws.Range("E38").Value = "2026-02-16 17:20:05"
$ws.Range("K38").Value = "13.8 MJ/m2"
$ws.Range("O38").Value = "11.7 °C"
$ws.Range("E39").Value = "2026-02-16 17:20:08"
$ws.Range("E40").Value = "2026-02-16 17:20:10"
$ws.Range("J40").Value = "1016.8 hPa"
$ws.Range("O40").Value = "6.6 °C"
$ws.Range("E41").Value = "2026-02-16 17:20:13"
$ws.Range("J41").Value = "1015.1 hPa"
$ws.Range("K41").Value = "11.2 MJ/m2"
$ws.Range("O41").Value = "17.3 °C"
$ws.Range("E42").Value = "2026-02-16 17:20:15"
$ws.Range("O42").Value = "11.0 °C"
$ws.Range("E43").Value = "2026-02-16 17:20:17"
$styleSource.Copy() | Out-Null
$ws.Range("H43").NumberFormat = "@"
$ws.Range("H43").Value = "81%"
$ws.Range("H43").PasteSpecial(-4122) | Out-Null
$ws.Range("K43").Value = "12.9 MJ/m2"
$ws.Range("O43").Value = "7.9 °C"
$ws.Range("E44").Value = "2026-02-16 17:20:19"
$ws.Range("K44").Value = "8.4 MJ/m2"
$ws.Range("N44").Value = "-1.3 °C 16:59 TU"
$ws.Range("E45").Value = "2026-02-16 17:20:22"
$ws.Range("I45").Value = "14.4 mm"
$ws.Range("K45").Value = "4.1 MJ/m2"
$ws.Range("E46").Value = "2026-02-16 17:20:24"
$styleSource.Copy() | Out-Null
$ws.Range("H46").NumberFormat = "@"
$ws.Range("H46").Value = "55%"
$ws.Range("H46").PasteSpecial(-4122) | Out-Null
$ws.Range("J46").Value = "1017.2 hPa"
$ws.Range("K46").Value = "12.7 MJ/m2"
$ws.Range("O46").Value = "15.8 °C"

$excel.CutCopyMode = 0
